$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ME-778"
$ws.Range("B2").Value = "OF-783"
$ws.Range("C2").Value = "CT-391"
$ws.Range("D2").Value = "BR-573"
$ws.Range("E2").Value = "Please"
$ws.Range("F2").Value = "RE-147"
$ws.Range("G2").Value = "OR-275"
$ws.Range("I2").Value = 9

$ws.Range("I6").Select()
